$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.139.54"
$ws.Range("E2").Value = "  -3.49%  "
$ws.Range("D3").Value = "1.603.19"
$ws.Range("E3").Value = "  -2.91%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.61"
$ws.Range("E6").Value = "  -2.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3782"
$ws.Range("E7").Value = "  -3.07%  "
$ws.Range("E8").Value = "  -4.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.06"
$ws.Range("E9").Value = "  -4.48%  "
$ws.Range("E10").Value = "  -6.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08157"
$ws.Range("E11").Value = "  -3.40%  "
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.06"
$ws.Range("E13").Value = "  -3.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.638"
$ws.Range("E14").Value = "  -6.18%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.438"
$ws.Range("E15").Value = "  -7.24%  "
$ws.Range("E16").Value = "  -4.09%  "
$ws.Range("D17").Value = "1.612.46"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.60"
$ws.Range("E18").Value = "  -3.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06843"
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.29"
$ws.Range("E20").Value = "  -6.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.589"
$ws.Range("E21").Value = "  -5.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.06"
$ws.Range("E23").Value = "  -5.18%  "
$ws.Range("D24").Value = "23.151.65"
$ws.Range("E24").Value = "  -3.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.343"
$ws.Range("E25").Value = "  -4.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.820"
$ws.Range("E26").Value = "  -4.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.10"
$ws.Range("E27").Value = "  -4.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.76"
$ws.Range("E28").Value = "  -1.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.285"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "131.98"
$ws.Range("E30").Value = "  -4.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.414"
$ws.Range("E31").Value = "  -3.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.898"
$ws.Range("E32").Value = "  -12.98%  "
$ws.Range("D33").Value = "1.787.18"
$ws.Range("E33").Value = "  -2.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07720"
$ws.Range("E34").Value = "  -4.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9478"
$ws.Range("E35").Value = "  -7.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02778"
$ws.Range("E36").Value = "  -5.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.280"
$ws.Range("E37").Value = "  -6.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2547"
$ws.Range("E38").Value = "  -4.82%  "
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.08913"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.11"
$ws.Range("E40").Value = "  -5.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.389"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.77"
$ws.Range("E42").Value = "  -4.71%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7117"
$ws.Range("E43").Value = "  -6.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.54"
$ws.Range("E44").Value = "  -4.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6638"
$ws.Range("E45").Value = "  -4.64%  "
$ws.Range("E46").Value = "  +0.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.306"
$ws.Range("E47").Value = "  -6.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.977"
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "131.81"
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("E50").Value = "  -4.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.212"
$ws.Range("E51").Value = "  -1.02%  "
